$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: LinearRegression - only B2 changes
$ws.Range("B2").Value = 3227498297016498

# Row 3: RandomForestRegressor - B3, C3, D3 change
$ws.Range("B3").Value = 1425813656586818
$ws.Range("C3").Value = 1432977807671950
$ws.Range("D3").Value = 2474058268104243

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, B4, C4, D4 change
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 32429314004701.44
$ws.Range("C4").Value = 32525312907575.65
$ws.Range("D4").Value = 2975661411605120

# Row 5: AdaBoostRegressor -> MLPRegressor, B5, C5, D5 change
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 508961369873299.8
$ws.Range("C5").Value = 248036081681645.1
$ws.Range("D5").Value = 456570472770661.2
